$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "W:\Data\Forecast\Tools\forecast_git\create_forecast_basic\current"
$ws.Range("B3").Value = "W:\Data\Forecast\forecast_by_version\V4\BASE_YEAR"
$ws.Range("B4").Value = "W:\Data\Forecast\forecast_by_version\V4\BASE_YEAR"
